$d = $word.ActiveDocument

# 1. Remove the old "_GoBack" bookmark that currently sits after
#    "... schneller aus." (it will be re-created further down, right
#    before the "dahbur2011survey" citation).
$oldBookmark = $d.Bookmarks("_GoBack")
$oldBookmark.Delete()

# 2. Locate "Beispiel BluePill [" and collapse the found range to its end
#    so we can insert the new citation text right after it.
$find = $d.Content
$found = $find.Find.Execute("Beispiel BluePill [", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$find.Collapse(0)
$insertPos = $find.Start

# 3. Insert the new paper reference text: "rutkowska2008bluepilling, "
$find.InsertAfter("rutkowska2008bluepilling, ")
$afterInsertPos = $insertPos + "rutkowska2008bluepilling, ".Length

# 4. Re-create the "_GoBack" bookmark right after the newly inserted text,
#    immediately before the existing "dahbur2011survey" citation.
$bookmarkRange = $d.Range($afterInsertPos, $afterInsertPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
